$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh: update Price (D) and Volume/1h (E) columns,
# and fix the Uniswap/Polkadot row order swap (rows 19-20).

# Row 2
$ws.Range("D2").Value = '70.581.05'
$ws.Range("E2").Value = '  +2.11%  '

# Row 3
$ws.Range("D3").Value = '3.804.00'
$ws.Range("E3").Value = '  +0.73%  '

# Row 4
$ws.Range("E4").Value = '  -0.04%  '

# Row 5
$ws.Range("D5").Value = '667.29'
$ws.Range("E5").Value = '  +6.73%  '

# Row 6
$ws.Range("D6").Value = '169.16'
$ws.Range("E6").Value = '  +1.93%  '

# Row 7
$ws.Range("D7").Value = '3.802.26'
$ws.Range("E7").Value = '  +0.78%  '

# Row 8
$ws.Range("E8").Value = '  -0.01%  '

# Row 9
$ws.Range("D9").Value = '0.525'
$ws.Range("E9").Value = '  +1.13%  '

# Row 10
$ws.Range("E10").Value = '  +0.30%  '

# Row 11
$ws.Range("E11").Value = '  +1.43%  '

# Row 12
$ws.Range("D12").Value = '7.01'
$ws.Range("E12").Value = '  +4.45%  '

# Row 13
$ws.Range("D13").Value = '0.0000244'
$ws.Range("E13").Value = '  -0.77%  '

# Row 14
$ws.Range("D14").Value = '35.79'
$ws.Range("E14").Value = '  +0.31%  '

# Row 15
$ws.Range("D15").Value = '4.443.99'
$ws.Range("E15").Value = '  +0.67%  '

# Row 16
$ws.Range("D16").Value = '3.798.22'
$ws.Range("E16").Value = '  +0.82%  '

# Row 17
$ws.Range("D17").Value = '70.533.45'
$ws.Range("E17").Value = '  +2.00%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.70'
$ws.Range("E18").Value = '  +0.12%  '

# Row 19
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = '7.19'
$ws.Range("E19").Value = '  +1.15%  '

# Row 20
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '11.61'
$ws.Range("E20").Value = '  +20.71%  '

# Row 21
$ws.Range("E21").Value = '  +0.56%  '

# Row 22
$ws.Range("D22").Value = '474.91'
$ws.Range("E22").Value = '  +1.45%  '

# Row 23
$ws.Range("D23").Value = '0.715'
$ws.Range("E23").Value = '  +1.25%  '

# Row 24
$ws.Range("D24").Value = '83.04'
$ws.Range("E24").Value = '  -0.10%  '

# Row 25
$ws.Range("D25").Value = '0.0000144'
$ws.Range("E25").Value = '  -2.49%  '

# Row 26
$ws.Range("D26").Value = '12.21'
$ws.Range("E26").Value = '  +1.69%  '

# Row 27
$ws.Range("D27").Value = '10.36'
$ws.Range("E27").Value = '  +3.42%  '

# Row 28
$ws.Range("E28").Value = '  -1.71%  '

# Row 29
$ws.Range("E29").Value = '  +0.06%  '

# Row 30
$ws.Range("D30").Value = '3.956.78'
$ws.Range("E30").Value = '  +0.75%  '

# Row 31
$ws.Range("D31").Value = '2.86'
$ws.Range("E31").Value = '  +6.83%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.30'
$ws.Range("E32").Value = '  +2.68%  '

# Row 33
$ws.Range("D33").Value = '7.39'
$ws.Range("E33").Value = '  +1.95%  '

# Row 34
$ws.Range("D34").Value = '29.63'
$ws.Range("E34").Value = '  +2.98%  '

# Row 35
$ws.Range("D35").Value = '0.177'
$ws.Range("E35").Value = '  +6.72%  '

# Row 36
$ws.Range("D36").Value = '9.12'
$ws.Range("E36").Value = '  +1.43%  '

# Row 37
$ws.Range("E37").Value = '  -0.03%  '

# Row 38
$ws.Range("D38").Value = '3.761.16'
$ws.Range("E38").Value = '  +0.86%  '

# Row 39
$ws.Range("E39").Value = '  +0.64%  '

# Row 40
$ws.Range("D40").Value = '3.42'
$ws.Range("E40").Value = '  +0.34%  '

# Row 41
$ws.Range("D41").Value = '5.97'
$ws.Range("E41").Value = '  +2.86%  '

# Row 42
$ws.Range("D42").Value = '0.967'
$ws.Range("E42").Value = '  +0.08%  '

# Row 43
$ws.Range("E43").Value = '  +0.05%  '

# Row 44
$ws.Range("D44").Value = '2.11'
$ws.Range("E44").Value = '  +9.93%  '

# Row 46
$ws.Range("D46").Value = '45.74'
$ws.Range("E46").Value = '  +5.86%  '

# Row 47
$ws.Range("D47").Value = '158.67'
$ws.Range("E47").Value = '  +4.43%  '

# Row 48
$ws.Range("E48").Value = '  +2.73%  '

# Row 49
$ws.Range("E49").Value = '  +4.95%  '

# Row 50
$ws.Range("E50").Value = '  +0.75%  '

# Row 51
$ws.Range("D51").Value = '8.52'
$ws.Range("E51").Value = '  +1.30%  '
